$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-53 (date serial, y_0, y_0_forecast, y_1, y_1_forecast)
$data = @(
    @(39400, 2007, 1.144978573787081, 2008, 2.918510996763723),
    @(39583, 2008, 3.09825757489699, 2009, 7.865470614547343),
    @(39765, 2008, 1.381024225294869, 2009, 4.264380119800992),
    @(39948, 2009, -1.791203563722299, 2010, -6.760862998203643),
    @(40130, 2009, -0.8792832172735965, 2010, 3.109784435759599),
    @(40310, 2010, 0.5799958470386946, 2011, 5.643342995751777),
    @(40494, 2010, 0.9337833426867448, 2011, 2.730731696345146),
    @(40676, 2011, 3.181454202131073, 2012, 4.38978860149748),
    @(40862, 2011, 2.791140000794279, 2012, 1.683857142130885),
    @(41044, 2012, 0.5930547804883668, 2013, -1.194610791900008),
    @(41228, 2012, 0.4451370000809973, 2013, -0.2532347529486723),
    @(41409, 2013, -0.3951783438669754, 2014, 0.03694906323863378),
    @(41592, 2013, 0.2545814083968478, 2014, 2.553470871380537),
    @(41774, 2014, 3.292216014290039, 2015, 7.617133650412211),
    @(41957, 2014, 1.297015177357297, 2015, -0.434146007584113),
    @(42137, 2015, 1.670328650030184, 2016, 2.037906845818616),
    @(42321, 2015, 1.365576377841027, 2016, 2.383242923544526),
    @(42503, 2016, 2.562791874943371, 2017, 3.265947405805814),
    @(42689, 2016, 2.204449574611278, 2017, 1.688977015142101),
    @(42867, 2017, 1.526411006965533, 2018, 0.6601843988560674),
    @(43053, 2017, 2.18621550610123, 2018, 2.066615940231964),
    @(43145, 2018, 1.297923389414657, 2019, 0.5447775838346658),
    @(43235, 2018, 1.63465618619294, 2019, 1.551857746372698),
    @(43326, 2018, 2.054458927584024, 2019, 3.441981941009353),
    @(43418, 2018, 1.911050033324102, 2019, 3.0862758122153),
    @(43510, 2019, 1.064009474888983, 2020, 0.03490120525229123),
    @(43600, 2019, 1.35261353265177, 2020, 0.8024032016000104),
    @(43691, 2019, 1.668617211002466, 2020, 1.816757311461803),
    @(43783, 2019, 1.457852003181337, 2020, -1.135072001636328),
    @(43875, 2020, 1.60064760240497, 2021, 2.904532120297287),
    @(43966, 2020, -2.082763426755907, 2021, -5.866344937500023),
    @(44068, 2020, -3.357986809108382, 2021, 11.62806235225531),
    @(44159, 2020, -3.258619210312885, 2021, -2.878617960200258),
    @(44251, 2021, -3.88825249955117, 2022, -7.923811177410267),
    @(44341, 2021, -0.1380317107957718, 2022, 7.749494937649115),
    @(44432, 2021, -0.1388955462784724, 2022, 4.532186626383039),
    @(44525, 2021, 0.4255262881966981, 2022, 1.466936654457096),
    @(44617, 2022, 3.299288015397095, 2023, 3.813885712818554),
    @(44706, 2022, 3.848999231984762, 2023, 2.866003071127765),
    @(44798, 2022, 3.463320568938566, 2023, 0.9167463358189964),
    @(44890, 2022, 3.293290997728171, 2023, -1.421977974472588),
    @(44981, 2023, -0.4618455958399603, 2024, -2.262006095280478),
    @(45071, 2023, -0.3745803349312071, 2024, 0.645722451525943),
    @(45163, 2023, -0.2871245688614854, 2024, -0.1008213866759977),
    @(45254, 2023, -0.2814561130375925, 2024, -0.3873858053678236),
    @(45345, 2024, -0.4631846496550684, 2025, 0.11654623407098),
    @(45436, 2024, -1.091476630333243, 2025, -0.4907904687545206),
    @(45534, 2024, -0.9636841177852018, 2025, -0.1102571493959759),
    @(45618, 2024, -0.6470065423293869, 2025, 1.276847713071927),
    @(45713, 2025, 2.749014873207956, 2026, 1.993341940459148),
    @(45800, 2025, 1.626992717807862, 2026, 0.6270138473519316),
    @(45891, 2025, 2.053865394798304, 2026, 0.8600861498751833),
)

# First extend the date-format style (currently applied through row 52 in column A)
# down to the newly added row 53 by copying the format from A52.
$ws.Range("A52").Copy()
$ws.Range("A53").PasteSpecial(-4122)
$excel.CutCopyMode = 0

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
